# Kientrucdetai_26_8.pptx - "Chinh sua ngay 31/8"
#
# 1) The auto-updating "datetimeFigureOut" field shown on every slide
#    master / slide layout footer was refreshed from 8/29/2010 to
#    8/30/2010 (the author re-saved the deck a day later).
# 2) The "Flowchart: Magnetic Disk 26" shape on slide 1 was nudged down
#    (Top 234pt -> 246pt, i.e. y = 2971800 EMU -> 3124200 EMU).

$p = $ppt.ActivePresentation

$oldDate = "8/29/2010"
$newDate = "8/30/2010"

# --- helper: walk a Shapes collection and refresh the cached date text
function Update-DateShapes {
    param($shapes)

    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# --- 1) Slide master date placeholder
Update-DateShapes $p.SlideMaster.Shapes

# --- 1) Every slide layout's date placeholder
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DateShapes $layouts.Item($i).Shapes
}

# --- 2) Move the "Flowchart: Magnetic Disk 26" shape on slide 1 down
$slide1 = $p.Slides.Item(1)
for ($k = 1; $k -le $slide1.Shapes.Count; $k++) {
    $shape = $slide1.Shapes.Item($k)
    if ($shape.Name -eq "Flowchart: Magnetic Disk 26") {
        $shape.Top = 246
    }
}
